# Apply corrected table values after realizing outputs were off from
# running only one section of the script.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new OTU/family label, new median_rank, new median_feature_weight
# ($null entries mean that column's value is unchanged for that row)
$updates = @(
    @(1, 'OTU', $null, $null),
    @(2, 'Desulfovibrio (OTU 3465)', 2.5, -0.3296925030340111),
    @(3, 'Ruminococcaceae (OTU 467)', 3.5, 0.4074428019338807),
    @(4, 'Lachnospiraceae (OTU 543)', 4, -0.3400803900547354),
    @(5, 'Erysipelotrichaceae (OTU 189)', 7, 0.3535074191234203),
    @(6, 'Porphyromonadaceae (OTU 228)', 10.5, 0.3180013784902163),
    @(7, 'Betaproteobacteria (OTU 58)', 13, -0.1839360244793108),
    @(8, 'Enterobacteriaceae (OTU 1)', 14.5, -0.1838742192593545),
    @(9, 'Enterococcus (OTU 23)', 17, -0.1692471862407768),
    @(10, 'Porphyromonadaceae (OTU 139)', 19.5, 0.2965338002606991),
    @(11, 'Lachnospiraceae (OTU 397)', 19.5, -0.2597015998577158),
    @(12, 'Burkholderiales (OTU 34)', 23, -0.1733591678659266),
    @(13, 'Porphyromonadaceae (OTU 87)', 23.5, 0.2699772801286908),
    @(14, 'Clostridium (OTU 154)', 26, -0.1562934472421492),
    @(15, 'Porphyromonadaceae (OTU 608)', 26, 0.2711412609900756),
    @(16, 'Porphyromonadaceae (OTU 222)', 30, 0.253072086268029),
    @(17, 'Ruminococcaceae (OTU 520)', 30.5, -0.1553644695555795),
    @(18, 'Bacillus (OTU 636)', 31.5, -0.1553644695555795),
    @(19, 'Dactylosporangium (OTU 3207)', 32.5, -0.1553644695555795),
    @(20, 'Coriobacteriaceae (OTU 293)', 34.5, 0.2492104073905196),
    @(21, 'Clostridiales (OTU 356)', 36.5, -0.243039217851724),
    @(22, 'Enterobacteriaceae (OTU 1)', $null, 0.4805050429743277),
    @(23, 'Bacteroides (OTU 2)', 3, -0.452390011750085),
    @(24, 'Lactobacillus (OTU 18)', 4, -0.4092711501860915),
    @(25, 'Escherichia/Shigella (OTU 610)', 6.5, -0.3211501162625571),
    @(26, 'Lachnospiraceae (OTU 56)', 13, 0.2834676535132536),
    @(27, 'Ruminococcaceae (OTU 520)', 13.5, -0.287966761885177),
    @(28, 'Porphyromonadaceae (OTU 54)', 15, -0.2011592266132916),
    @(29, 'Porphyromonadaceae (OTU 22)', 16, 0.2237962690648616),
    @(30, 'Lachnospiraceae (OTU 38)', 16, -0.08215621458561985),
    @(31, 'Lachnospiraceae (OTU 33)', 17, -0.2524117226522107),
    @(32, 'Porphyromonadaceae (OTU 7)', 20, 0.1314449678811958),
    @(33, 'Ruminococcaceae (OTU 60)', 24, -0.1139403081751913),
    @(34, 'Erysipelotrichaceae (OTU 234)', 24.5, 0.1968104657738282),
    @(35, 'Lachnospiraceae (OTU 9)', 27, -0.2105210354545368),
    @(36, 'Proteus (OTU 16)', 27, -0.190954390712031),
    @(37, 'Alishewanella (OTU 776)', 27, -0.1122260698812191),
    @(38, 'Eisenbergiella (OTU 164)', 32.5, -0.0456374695110339),
    @(39, 'Clostridium (OTU 99)', 33, 0.1805555754102407),
    @(40, 'Clostridium (OTU 226)', 33, 0.04605844444543253),
    @(41, 'Lactobacillus (OTU 834)', 33.5, -0.0456374695110339),
    @(42, 'Lactobacillus (OTU 18)', $null, -1.209561846396649),
    @(43, 'Bacteroides (OTU 2)', $null, -1.044768312738835),
    @(44, 'Lachnospiraceae (OTU 35)', $null, 0.8395257141727522),
    @(45, 'Coriobacteriaceae (OTU 3419)', 9, -0.4249008232572614),
    @(46, 'Coriobacteriaceae (OTU 379)', 10.5, 0.8030734725326898),
    @(47, 'Turicibacter (OTU 14)', 11.5, 0.6296225168141147),
    @(48, 'Bifidobacterium (OTU 46)', $null, 0.7453827213448163),
    @(49, 'Bacteroides (OTU 3)', 12.5, 0.6747382360394445),
    @(50, 'Enterococcus (OTU 23)', 13.5, -0.6034903369986337),
    @(51, 'Porphyromonadaceae (OTU 7)', 17.5, 0.6152684912951247),
    @(52, 'Erysipelotrichaceae (OTU 234)', 17.5, 0.6021268407672054),
    @(53, 'Bacteria (OTU 509)', 18, 0.6066825319705877),
    @(54, 'Lachnospiraceae (OTU 44)', 20, 0.5807993220192496),
    @(55, 'Lactobacillus (OTU 31)', 21, 0.5531318973366881),
    @(56, 'Lactobacillales (OTU 603)', 21, 0.5721851142342991),
    @(57, 'Bacteroides (OTU 32)', 23, 0.1162382837427073),
    @(58, 'Alistipes (OTU 161)', 26, 0.09003603176770333),
    @(59, 'Anaerofustis (OTU 475)', 26, -0.4718461291268826),
    @(60, 'Lachnospiraceae (OTU 109)', 28.5, 0.009322759602494523),
    @(61, 'Streptococcus (OTU 512)', 28.5, 0.530449090825015)
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 1).Value = $u[1]
    if ($null -ne $u[2]) {
        $ws.Cells.Item($r, 2).Value = $u[2]
    }
    if ($null -ne $u[3]) {
        $ws.Cells.Item($r, 3).Value = $u[3]
    }
}
